# Apply cryptocurrency price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text interpretation, matching the original
# inline-string cell content (avoids Excel auto-converting values such as
# "9.60" or "1.00" into numbers and stripping trailing zeros).

$ws.Range('D2').Value = "'26.722.34"
$ws.Range('E2').Value = "'  +0.11%  "
$ws.Range('D3').Value = "'1.647.98"
$ws.Range('E3').Value = "'  +0.73%  "
$ws.Range('E4').Value = "'  +0.25%  "
$ws.Range('D5').Value = "'216.21"
$ws.Range('E5').Value = "'  +1.38%  "
$ws.Range('D6').Value = "'0.505"
$ws.Range('E6').Value = "'  -0.73%  "
$ws.Range('E7').Value = "'  +0.24%  "
$ws.Range('E8').Value = "'  -0.30%  "
$ws.Range('D9').Value = "'0.0627"
$ws.Range('E9').Value = "'  +0.71%  "
$ws.Range('D10').Value = "'19.45"
$ws.Range('E10').Value = "'  +1.11%  "
$ws.Range('D11').Value = "'0.0845"
$ws.Range('E11').Value = "'  +0.31%  "
$ws.Range('D12').Value = "'1.878.63"
$ws.Range('B13').Value = "'Polkadot"
$ws.Range('C13').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D13').Value = "'4.24"
$ws.Range('E13').Value = "'  +3.38%  "
$ws.Range('B14').Value = "'WrappedEther"
$ws.Range('C14').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D14').Value = "'1.643.13"
$ws.Range('E14').Value = "'  +0.53%  "
$ws.Range('E15').Value = "'  +1.54%  "
$ws.Range('D16').Value = "'66.38"
$ws.Range('E16').Value = "'  +4.87%  "
$ws.Range('D17').Value = "'26.757.21"
$ws.Range('E17').Value = "'  +0.30%  "
$ws.Range('D18').Value = "'0.0₃0757"
$ws.Range('E18').Value = "'  +1.42%  "
$ws.Range('D19').Value = "'220.29"
$ws.Range('E19').Value = "'  +0.89%  "
$ws.Range('E20').Value = "'  +0.28%  "
$ws.Range('D21').Value = "'4.41"
$ws.Range('E21').Value = "'  +2.11%  "
$ws.Range('E22').Value = "'  +2.13%  "
$ws.Range('D23').Value = "'9.60"
$ws.Range('E23').Value = "'  +1.65%  "
$ws.Range('D24').Value = "'2.10"
$ws.Range('E24').Value = "'  +9.21%  "
$ws.Range('E25').Value = "'  -0.74%  "
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = "'  +0.20%  "
$ws.Range('E27').Value = "'  -0.67%  "
$ws.Range('E28').Value = "'  +2.71%  "
$ws.Range('D29').Value = "'15.93"
$ws.Range('E29').Value = "'  +2.77%  "
$ws.Range('E30').Value = "'  +1.86%  "
$ws.Range('E31').Value = "'  +0.76%  "
$ws.Range('E32').Value = "'  +3.15%  "
$ws.Range('D33').Value = "'3.08"
$ws.Range('E33').Value = "'  +2.97%  "
$ws.Range('D34').Value = "'1.289.22"
$ws.Range('E34').Value = "'  +7.53%  "
$ws.Range('D35').Value = "'1.56"
$ws.Range('E35').Value = "'  +2.86%  "
$ws.Range('D36').Value = "'0.0186"
$ws.Range('D37').Value = "'2.40"
$ws.Range('D38').Value = "'0.834"
$ws.Range('E38').Value = "'  +2.80%  "
$ws.Range('D39').Value = "'0.527"
$ws.Range('E39').Value = "'  +4.29%  "
$ws.Range('E40').Value = "'  +0.27%  "
$ws.Range('D41').Value = "'0.812"
$ws.Range('E41').Value = "'  +2.23%  "
$ws.Range('E42').Value = "'  -1.72%  "
$ws.Range('D43').Value = "'5.45"
$ws.Range('E43').Value = "'  +0.48%  "
$ws.Range('D44').Value = "'1.789.39"
$ws.Range('E44').Value = "'  +1.13%  "
$ws.Range('D45').Value = "'93.88"
$ws.Range('E45').Value = "'  +1.40%  "
$ws.Range('D46').Value = "'60.73"
$ws.Range('E46').Value = "'  +10.71%  "
$ws.Range('D47').Value = "'1.62"
$ws.Range('E47').Value = "'  +3.93%  "
$ws.Range('E48').Value = "'  +0.81%  "
$ws.Range('D49').Value = "'7.85"
$ws.Range('E49').Value = "'  +2.59%  "
$ws.Range('D50').Value = "'0.0980"
$ws.Range('E50').Value = "'  +3.21%  "
$ws.Range('E51').Value = "'  -0.66%  "
